$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 12 with "Testing" text in B12 (appends a new shared string
# "Testing" and extends the sheet's used range / dimension to B2:F12)
$ws.Range("B12").Value = "Testing"

# Update the selection to match the newly entered cell (was B13, now B12)
$ws.Range("B12").Select()
